$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.650.90"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "1.812.95"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.600"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.62%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "38.04"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.76%  "
$ws.Range("E9").Value = "  -3.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0681"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0972"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.61%  "
$ws.Range("D12").Value = "2.073.78"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("D14").Value = "1.807.05"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.635"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.32%  "
$ws.Range("D16").Value = "34.604.75"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.45"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.80"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("E20").Value = "  -2.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.52%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("E24").Value = "  +4.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.54%  "
$ws.Range("E28").Value = "  +1.86%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.38%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.35%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0525"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.30%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.05%  "
$ws.Range("E34").Value = "  -0.54%  "
$ws.Range("D35").Value = "1.366.83"
$ws.Range("E35").Value = "  -2.29%  "
$ws.Range("E36").Value = "  -3.35%  "
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("E38").Value = "  -1.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.35"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.33%  "
$ws.Range("E40").Value = "  +6.88%  "
$ws.Range("E41").Value = "  +1.48%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.54%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "81.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.19%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.942"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.42%  "
$ws.Range("E46").Value = "  -2.19%  "
$ws.Range("D47").Value = "1.974.82"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("E48").Value = "  -3.37%  "
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.40%  "
$ws.Range("D51").Value = "0.0₆0121"
$ws.Range("E51").Value = "  -7.51%  "
